$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: "Datos Exportados" -> "Producto"
$wb.Sheets.Item(1).Name = "Producto"

# Drop the old stock/last-update columns (C, D) entirely: removes their
# column-width defs and their cell content while leaving the A/B row
# striping styles (s=1/2/3) untouched on the remaining cells.
$ws.Range("C1:D8").EntireColumn.Delete()

# Wipe the old A/B content (keeps per-row styles intact) before writing
# the new "Producto" characteristics sheet.
$ws.Range("A1:B8").ClearContents()

# New column widths: A=30 chars, B=50 chars. The COM bridge rounds
# ColumnWidth to a pixel grid with ~0.83-char padding baked in, so asking
# for X.09 lands exactly on the integral stored width of X+1.
$ws.Columns.Item(1).ColumnWidth = 29.09
$ws.Columns.Item(2).ColumnWidth = 49.09

# Section title (merged across A1:B1)
$ws.Range("A1").Value = "CARACTERISTICAS GENERALES"

$ws.Range("A2").Value = "NOMBRE"
$ws.Range("B2").Value = "Cerveza Artesanal Actualizada"

$ws.Range("A3").Value = "DESCRIPCIÓN"
$ws.Range("B3").Value = "Cerveza artesanal de alta calidad, edición limitada"

$ws.Range("A4").Value = "MARCA"
$ws.Range("B4").Value = "Artesanal Premium"

$ws.Range("A5").Value = "CONTENIDO"
$ws.Range("B5").Value = "'30"

$ws.Range("A6").Value = "UNIDAD DE MEDIDA"
$ws.Range("B6").Value = "L"

$ws.Range("A7").Value = "TIPO"
$ws.Range("B7").Value = "bebida"

$ws.Range("A8").Value = "PRECIO"
$ws.Range("B8").Value = "'1800"

$ws.Range("A9").Value = "CATEGORÍA"
$ws.Range("B9").Value = "Bebidas alcohólicas"

# Row 9 is brand new (sheet used to stop at row 8) so it has no style yet;
# clone the alternating-fill look from row 7 (same odd-row style) without
# minting a new style entry.
$ws.Range("A7:B7").Copy()
$ws.Range("A9:B9").PasteSpecial(-4122)

# Title bar spans both columns
$ws.Range("A1:B1").Merge()
